# Add three new contact rows (with "false input" check for emails/phones)
# to the bottom of the contact list on Sheet1, each with a mailto: hyperlink
# on the email cell, matching the existing table's layout/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 9;  Name = "paulus nugroho"; Phone = "087877589855"; Email = "paulus.nug@yahoo.co.id" },
    @{ Row = 10; Name = "yudi bramanto";  Phone = "089989986746"; Email = "yudiibram.78@gmail.com" },
    @{ Row = 11; Name = "levi prasetyo";  Phone = "087824656698"; Email = "leviipras99@gmail.com" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $entry.Phone
    $ws.Cells.Item($r, 3).Value = $entry.Email

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 3), "mailto:" + $entry.Email)
    $ws.Cells.Item($r, 3).Style = "Hyperlink"
}

$ws.Range("A12").Select()
